$wb = $excel.ActiveWorkbook

# --- 1. Clean up the "Информация" sheet: delete stale / duplicate rows ---
# These are the rows that were previously hidden (filtered out, left-over
# duplicates from an old import) - they get physically removed from the sheet.
$infoSheet = $wb.Worksheets.Item("Информация")

$rowsToDelete = @(2, 6, 7, 11, 14, 17, 24, 25, 30, 31, 32, 37, 38, 43, 48, 52, 55, 58, 59, 63, 64, 65, 66, 73)
$rowsToDelete = $rowsToDelete | Sort-Object -Descending

foreach ($r in $rowsToDelete) {
    $infoSheet.Rows.Item($r).Delete()
}

# --- 2. Fix up an address that was missing its postal code ---
$found = $infoSheet.Cells.Item(31, 6)
$found.Value = "Sunny Pioneer Green, Wamduska, MI, 15786-3482"

# --- 3. Make "Информация" the active sheet/tab ---
$infoSheet.Activate()
